# Update the build timestamp embedded in the "version" text throughout the
# workbook, from "January 30 2026 16.19.47 EST" to
# "February 02 2026 12.49.33 EST".

$wb = $excel.ActiveWorkbook

$oldStamp = "January 30 2026 16.19.47 EST"
$newStamp = "February 02 2026 12.49.33 EST"

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        $val = $cell.Text
        if ($val -ne $null -and $val -is [string] -and $val.Contains($oldStamp)) {
            $cell.Value = $val.Replace($oldStamp, $newStamp)
        }
    }
}
